$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (D/E/F): state, language, theme ---
$ws.Range("D1").Value = "state"
$ws.Range("E1").Value = "language"
$ws.Range("F1").Value = "theme"

# --- Row 2 new data (written E, F, D to reproduce the shared-string insertion order of the target file) ---
$ws.Range("E2").Value = "English,Hindi-हिन्दी,Kannada-ಕೆನಡಾ,Malayalam-മലയാളം,Other Language,Tamil-தமிழ்,Telugu-తెలుగు"
$ws.Range("F2").Value = "Sustainable Development and Environment,Digital Transformation,Health and Well-being,Quality Education,Economic Empowerment,Smart and Resilient Communities,Agriculture and Rural Development,Others"
$ws.Range("D2").Value = "Andaman and Nicobar Islands,Andhra Pradesh,Arunachal Pradesh,Assam,Bihar,Chandigarh,Chhattisgarh,Dadra and Nagar Haveli and Daman and Diu,Delhi,Goa,Gujarat,Haryana,Himachal Pradesh,Jammu and Kashmir,Jharkhand,Karnataka,Kerala,Ladakh,Lakshadweep,Madhya Pradesh,Maharashtra,Manipur,Meghalaya,Mizoram,Nagaland,Odisha,Puducherry,Punjab,Rajasthan,Sikkim,Tamil Nadu,Telangana,Tripura,Uttar Pradesh,Uttarakhand,West Bengal"

# --- Row 3: B3 email changes, plus new D/E/F data ---
$ws.Range("B3").Value = "fbvcbvc@gmail.com"
$ws.Range("D3").Value = "Meghalaya,Mizoram,Nagaland,Sikkim,Tamil Nadu,Telangana,Tripura,Uttar Pradesh,Uttarakhand,West Bengal"
$ws.Range("E3").Value = "English"
$ws.Range("F3").Value = "Health and Well-being,Quality Education"

# --- Rows 4 & 5: drop the old record rows entirely (A/C cleared, B cleared but keeps its Hyperlink style) ---
$ws.Range("A4:C5").ClearContents()

# --- Hyperlinks: the engine's Hyperlinks collection only supports clearing
# the whole sheet at once, so drop them all and re-create the two that
# still exist in the target (B2, B3); B4/B5's are gone for good. ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:dzczxcvxz@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:fbvcbvc@gmail.com")
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"

# --- Selection moves to F3 ---
$ws.Range("F3").Select()
